# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers section (row 3 = Intel(R) Wi-Fi 6E AX211 160MHz - 23.60.1.2)
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 98.59999999999999

# Totals row (row 4)
$ws.Range("C4").Value = 1

# Good Drivers section - Total Samples counts
$ws.Range("B12").Value = 11140
$ws.Range("B13").Value = 14487
